# Weekly update of the "Puerro" (leek) price sheet.
# Three new daily observations are inserted into the historical table,
# pushing the existing rows down. The new rows land (in final row
# numbering) at rows 31, 33 and 45.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PuerroRow {
    param($Row, $Fecha, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $PrecioKg)

    $ws.Cells.Item($Row, 1).Value = 9
    $ws.Cells.Item($Row, 2).Value = "Vega Central Mapocho de Santiago"
    $ws.Cells.Item($Row, 3).Value = "Metropolitana"
    $ws.Cells.Item($Row, 4).Value = $Fecha
    $ws.Cells.Item($Row, 5).Value = 13
    $ws.Cells.Item($Row, 6).Value = 100112005
    $ws.Cells.Item($Row, 7).Value = "Puerro"
    $ws.Cells.Item($Row, 8).Value = "Sin especificar"
    $ws.Cells.Item($Row, 9).Value = "Primera"
    $ws.Cells.Item($Row, 10).Value = $Volumen
    $ws.Cells.Item($Row, 11).Value = $PrecioMin
    $ws.Cells.Item($Row, 12).Value = $PrecioMax
    $ws.Cells.Item($Row, 13).Value = $PrecioProm
    $ws.Cells.Item($Row, 14).Value = "`$/paquete 20 unidades"
    $ws.Cells.Item($Row, 15).Value = "Provincia de Chacabuco"
    $ws.Cells.Item($Row, 16).Value = $PrecioKg
    $ws.Cells.Item($Row, 17).Value = 20
    $ws.Cells.Item($Row, 18).Value = "Hortaliza"
}

# 1) Insert new row at 31 (existing rows 31-55 shift down to 32-56).
$ws.Rows(31).Insert()
Set-PuerroRow 31 44421 180 7000 8000 7500 375

# 2) Insert new row at 33 (old row 31, now at 32, stays put; rows from
#    33 on shift down again).
$ws.Rows(33).Insert()
Set-PuerroRow 33 44426 97 7000 8000 7505 375

# 3) Insert new row at 45 (old row 43, now at 44, stays put; rows from
#    45 on shift down again).
$ws.Rows(45).Insert()
Set-PuerroRow 45 44419 160 7000 8000 7500 375
